$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert two new columns before column D, shifting old D:K to F:M ---
$ws.Columns("D:E").Insert()

# --- Step 2: copy number/date formatting from column F (old column D) into new D:E ---
$ws.Columns("F").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 3: write the refreshed quarterly figures (columns D through M) ---
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42916
$ws.Range("K7").Value = 42825
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42643
$ws.Range("D8").Value = 1011800
$ws.Range("E8").Value = 940700
$ws.Range("F8").Value = 834100
$ws.Range("G8").Value = 825400
$ws.Range("H8").Value = 989500
$ws.Range("I8").Value = 939900
$ws.Range("J8").Value = 904300
$ws.Range("K8").Value = 897600
$ws.Range("L8").Value = 909900
$ws.Range("M8").Value = 992800
$ws.Range("D9").Value = 419600
$ws.Range("E9").Value = 500800
$ws.Range("F9").Value = 472200
$ws.Range("G9").Value = 482200
$ws.Range("H9").Value = 493700
$ws.Range("I9").Value = 555200
$ws.Range("J9").Value = 616900
$ws.Range("K9").Value = 553200
$ws.Range("L9").Value = 500500
$ws.Range("M9").Value = 556600
$ws.Range("D10").Value = 592200
$ws.Range("E10").Value = 439900
$ws.Range("F10").Value = 361900
$ws.Range("G10").Value = 343100
$ws.Range("H10").Value = 495900
$ws.Range("I10").Value = 384600
$ws.Range("J10").Value = 287400
$ws.Range("K10").Value = 344400
$ws.Range("L10").Value = 409400
$ws.Range("M10").Value = 436200
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = 3700
$ws.Range("E14").Value = 2500
$ws.Range("F14").Value = 500
$ws.Range("G14").Value = 1600
$ws.Range("H14").Value = 3600
$ws.Range("I14").Value = 3000
$ws.Range("J14").Value = 3100
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("D15").Value = 90000
$ws.Range("E15").Value = 85800
$ws.Range("F15").Value = 85800
$ws.Range("G15").Value = 54600
$ws.Range("H15").Value = 57800
$ws.Range("I15").Value = 55200
$ws.Range("J15").Value = 57600
$ws.Range("K15").Value = 57800
$ws.Range("L15").Value = 61400
$ws.Range("M15").Value = 61700
$ws.Range("D17").Value = 627100
$ws.Range("E17").Value = 696700
$ws.Range("F17").Value = 677800
$ws.Range("G17").Value = 624700
$ws.Range("H17").Value = 704500
$ws.Range("I17").Value = 708200
$ws.Range("J17").Value = 760600
$ws.Range("K17").Value = 702100
$ws.Range("L17").Value = 722800
$ws.Range("M17").Value = 768700
$ws.Range("D18").Value = 384800
$ws.Range("E18").Value = 244000
$ws.Range("F18").Value = 156300
$ws.Range("G18").Value = 200700
$ws.Range("H18").Value = 285000
$ws.Range("I18").Value = 231600
$ws.Range("J18").Value = 143700
$ws.Range("K18").Value = 195500
$ws.Range("L18").Value = 187100
$ws.Range("M18").Value = 224100
$ws.Range("D20").Value = -13300
$ws.Range("E20").Value = -20000
$ws.Range("F20").Value = -13600
$ws.Range("G20").Value = 7100
$ws.Range("H20").Value = 7200
$ws.Range("I20").Value = 4600
$ws.Range("J20").Value = 13200
$ws.Range("K20").Value = 168000
$ws.Range("L20").Value = 4800
$ws.Range("M20").Value = 184000
$ws.Range("D21").Value = 375700
$ws.Range("E21").Value = 224000
$ws.Range("F21").Value = 173900
$ws.Range("G21").Value = 262500
$ws.Range("H21").Value = 294700
$ws.Range("I21").Value = 236200
$ws.Range("J21").Value = 155700
$ws.Range("K21").Value = 421300
$ws.Range("L21").Value = 191600
$ws.Range("M21").Value = 407900
$ws.Range("D22").Value = 34400
$ws.Range("E22").Value = 33900
$ws.Range("F22").Value = 30900
$ws.Range("G22").Value = 14400
$ws.Range("H22").Value = 13400
$ws.Range("I22").Value = 14300
$ws.Range("J22").Value = 14700
$ws.Range("K22").Value = 16900
$ws.Range("L22").Value = 16700
$ws.Range("M22").Value = 16700
$ws.Range("D23").Value = 337000
$ws.Range("E23").Value = 190200
$ws.Range("F23").Value = 111900
$ws.Range("G23").Value = 193400
$ws.Range("H23").Value = 278700
$ws.Range("I23").Value = 221900
$ws.Range("J23").Value = 142200
$ws.Range("K23").Value = 346500
$ws.Range("L23").Value = 175200
$ws.Range("M23").Value = 391500
$ws.Range("D24").Value = 95200
$ws.Range("E24").Value = 54400
$ws.Range("F24").Value = 28300
$ws.Range("G24").Value = 47800
$ws.Range("H24").Value = 38600
$ws.Range("I24").Value = 55300
$ws.Range("J24").Value = 42500
$ws.Range("K24").Value = 76400
$ws.Range("L24").Value = 21700
$ws.Range("M24").Value = 83300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 241900
$ws.Range("E26").Value = 135800
$ws.Range("F26").Value = 83600
$ws.Range("G26").Value = 145600
$ws.Range("H26").Value = 240100
$ws.Range("I26").Value = 166600
$ws.Range("J26").Value = 99700
$ws.Range("K26").Value = 270200
$ws.Range("L26").Value = 153500
$ws.Range("M26").Value = 308100
$ws.Range("D27").Value = 225100
$ws.Range("E27").Value = 126500
$ws.Range("F27").Value = 77000
$ws.Range("G27").Value = 103100
$ws.Range("H27").Value = 149700
$ws.Range("I27").Value = 114500
$ws.Range("J27").Value = 78000
$ws.Range("K27").Value = 176100
$ws.Range("L27").Value = 104500
$ws.Range("M27").Value = 208900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = 13300
$ws.Range("E32").Value = 20000
$ws.Range("F32").Value = 13600
$ws.Range("G32").Value = -7100
$ws.Range("H32").Value = -7200
$ws.Range("I32").Value = -4600
$ws.Range("J32").Value = -13200
$ws.Range("K32").Value = -168000
$ws.Range("L32").Value = -4800
$ws.Range("M32").Value = -184000
$ws.Range("D33").Value = 225100
$ws.Range("E33").Value = 126500
$ws.Range("F33").Value = 77000
$ws.Range("G33").Value = 103100
$ws.Range("H33").Value = 149700
$ws.Range("I33").Value = 114500
$ws.Range("J33").Value = 78000
$ws.Range("K33").Value = 176100
$ws.Range("L33").Value = 104500
$ws.Range("M33").Value = 208900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("D35").Value = 225100
$ws.Range("E35").Value = 126500
$ws.Range("F35").Value = 77000
$ws.Range("G35").Value = 103100
$ws.Range("H35").Value = 149700
$ws.Range("I35").Value = 114500
$ws.Range("J35").Value = 78000
$ws.Range("K35").Value = 176100
$ws.Range("L35").Value = 104500
$ws.Range("M35").Value = 208900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42916
$ws.Range("K38").Value = 42825
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42643
$ws.Range("D41").Value = 360400
$ws.Range("E41").Value = 293100
$ws.Range("F41").Value = 150500
$ws.Range("G41").Value = 605900
$ws.Range("H41").Value = 616600
$ws.Range("I41").Value = 400000
$ws.Range("J41").Value = 189800
$ws.Range("K41").Value = 537400
$ws.Range("L41").Value = 371500
$ws.Range("M41").Value = 332600
$ws.Range("D42").Value = 1900
$ws.Range("E42").Value = 800
$ws.Range("F42").Value = 800
$ws.Range("G42").Value = 1433100
$ws.Range("H42").Value = 1100
$ws.Range("I42").Value = 400
$ws.Range("J42").Value = 300
$ws.Range("K42").Value = 700
$ws.Range("L42").Value = 700
$ws.Range("M42").Value = 700
$ws.Range("D43").Value = 929200
$ws.Range("E43").Value = 741200
$ws.Range("F43").Value = 820800
$ws.Range("G43").Value = 687700
$ws.Range("H43").Value = 818200
$ws.Range("I43").Value = 786200
$ws.Range("J43").Value = 793000
$ws.Range("K43").Value = 783300
$ws.Range("L43").Value = 835900
$ws.Range("M43").Value = 889400
$ws.Range("D44").Value = 83700
$ws.Range("E44").Value = 69600
$ws.Range("F44").Value = 63100
$ws.Range("G44").Value = 71500
$ws.Range("H44").Value = 58300
$ws.Range("I44").Value = 54100
$ws.Range("J44").Value = 55000
$ws.Range("K44").Value = 70500
$ws.Range("L44").Value = 56700
$ws.Range("M44").Value = 61500
$ws.Range("D45").Value = 90300
$ws.Range("E45").Value = 30500
$ws.Range("F45").Value = 33400
$ws.Range("G45").Value = 58000
$ws.Range("H45").Value = 63300
$ws.Range("I45").Value = 18800
$ws.Range("J45").Value = 26500
$ws.Range("K45").Value = 31200
$ws.Range("L45").Value = 43700
$ws.Range("M45").Value = 23300
$ws.Range("D46").Value = 1465500
$ws.Range("E46").Value = 1135200
$ws.Range("F46").Value = 1068500
$ws.Range("G46").Value = 2856300
$ws.Range("H46").Value = 1551100
$ws.Range("I46").Value = 1259400
$ws.Range("J46").Value = 1064700
$ws.Range("K46").Value = 1423100
$ws.Range("L46").Value = 1308500
$ws.Range("M46").Value = 1307600
$ws.Range("D47").Value = 112400
$ws.Range("E47").Value = 163200
$ws.Range("F47").Value = 163500
$ws.Range("G47").Value = 77500
$ws.Range("H47").Value = 85800
$ws.Range("I47").Value = 82700
$ws.Range("J47").Value = 81600
$ws.Range("K47").Value = 86200
$ws.Range("L47").Value = 83900
$ws.Range("M47").Value = 101600
$ws.Range("D48").Value = 7814800
$ws.Range("E48").Value = 7567400
$ws.Range("F48").Value = 7521100
$ws.Range("G48").Value = 5288300
$ws.Range("H48").Value = 5283200
$ws.Range("I48").Value = 5206300
$ws.Range("J48").Value = 5139800
$ws.Range("K48").Value = 5256500
$ws.Range("L48").Value = 5261200
$ws.Range("M48").Value = 5241700
$ws.Range("D49").Value = 1514700
$ws.Range("E49").Value = 1498200
$ws.Range("F49").Value = 1488500
$ws.Range("G49").Value = 1382400
$ws.Range("H49").Value = 1385400
$ws.Range("I49").Value = 1368600
$ws.Range("J49").Value = 1368000
$ws.Range("K49").Value = 1404500
$ws.Range("L49").Value = 1406900
$ws.Range("M49").Value = 1401800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("D52").Value = 99900
$ws.Range("E52").Value = 159200
$ws.Range("F52").Value = 140300
$ws.Range("G52").Value = 65500
$ws.Range("H52").Value = 69700
$ws.Range("I52").Value = 100800
$ws.Range("J52").Value = 97100
$ws.Range("K52").Value = 96700
$ws.Range("L52").Value = 91600
$ws.Range("M52").Value = 87100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = 11007400
$ws.Range("E54").Value = 10523300
$ws.Range("F54").Value = 10381900
$ws.Range("G54").Value = 9670000
$ws.Range("H54").Value = 8371300
$ws.Range("I54").Value = 8017800
$ws.Range("J54").Value = 7751100
$ws.Range("K54").Value = 8267000
$ws.Range("L54").Value = 8152100
$ws.Range("M54").Value = 8139700
$ws.Range("D57").Value = 331600
$ws.Range("E57").Value = 183700
$ws.Range("F57").Value = 209600
$ws.Range("G57").Value = 228600
$ws.Range("H57").Value = 312800
$ws.Range("I57").Value = 300700
$ws.Range("J57").Value = 277000
$ws.Range("K57").Value = 259900
$ws.Range("L57").Value = 283500
$ws.Range("M57").Value = 336000
$ws.Range("D58").Value = 484000
$ws.Range("E58").Value = 585500
$ws.Range("F58").Value = 145700
$ws.Range("G58").Value = 29900
$ws.Range("H58").Value = 25400
$ws.Range("I58").Value = 30000
$ws.Range("J58").Value = 26300
$ws.Range("K58").Value = 31500
$ws.Range("L58").Value = 27200
$ws.Range("M58").Value = 31200
$ws.Range("D59").Value = 974000
$ws.Range("E59").Value = 655700
$ws.Range("F59").Value = 601400
$ws.Range("G59").Value = 728700
$ws.Range("H59").Value = 900400
$ws.Range("I59").Value = 507400
$ws.Range("J59").Value = 439800
$ws.Range("K59").Value = 719000
$ws.Range("L59").Value = 832800
$ws.Range("M59").Value = 555900
$ws.Range("D60").Value = 1789700
$ws.Range("E60").Value = 1424900
$ws.Range("F60").Value = 956700
$ws.Range("G60").Value = 987200
$ws.Range("H60").Value = 1200700
$ws.Range("I60").Value = 838100
$ws.Range("J60").Value = 743100
$ws.Range("K60").Value = 1010400
$ws.Range("L60").Value = 1143400
$ws.Range("M60").Value = 923100
$ws.Range("D61").Value = 2503500
$ws.Range("E61").Value = 2423100
$ws.Range("F61").Value = 2950000
$ws.Range("G61").Value = 2476900
$ws.Range("H61").Value = 1118600
$ws.Range("I61").Value = 1145400
$ws.Range("J61").Value = 1174300
$ws.Range("K61").Value = 1206900
$ws.Range("L61").Value = 1211100
$ws.Range("M61").Value = 1344900
$ws.Range("D62").Value = 1313200
$ws.Range("E62").Value = 1287400
$ws.Range("F62").Value = 1234500
$ws.Range("G62").Value = 477500
$ws.Range("H62").Value = 485200
$ws.Range("I62").Value = 524400
$ws.Range("J62").Value = 538800
$ws.Range("K62").Value = 558500
$ws.Range("L62").Value = 568400
$ws.Range("M62").Value = 579500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("D66").Value = 5978200
$ws.Range("E66").Value = 5500500
$ws.Range("F66").Value = 5497900
$ws.Range("G66").Value = 5172300
$ws.Range("H66").Value = 3985700
$ws.Range("I66").Value = 3656200
$ws.Range("J66").Value = 3535100
$ws.Range("K66").Value = 3923000
$ws.Range("L66").Value = 3979300
$ws.Range("M66").Value = 3933200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("D72").Value = -544900
$ws.Range("E72").Value = -611100
$ws.Range("F72").Value = -742700
$ws.Range("G72").Value = 1245300
$ws.Range("H72").Value = 1146800
$ws.Range("I72").Value = 1154300
$ws.Range("J72").Value = 1037500
$ws.Range("K72").Value = 1081500
$ws.Range("L72").Value = 905400
$ws.Range("M72").Value = 944400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("D76").Value = 5029200
$ws.Range("E76").Value = 5022700
$ws.Range("F76").Value = 4884000
$ws.Range("G76").Value = 4497700
$ws.Range("H76").Value = 4385600
$ws.Range("I76").Value = 4361600
$ws.Range("J76").Value = 4216000
$ws.Range("K76").Value = 4343900
$ws.Range("L76").Value = 4172700
$ws.Range("M76").Value = 4206500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42916
$ws.Range("K80").Value = 42825
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42643
$ws.Range("D81").Value = 225100
$ws.Range("E81").Value = 126500
$ws.Range("F81").Value = 77000
$ws.Range("G81").Value = 103100
$ws.Range("H81").Value = 149700
$ws.Range("I81").Value = 114500
$ws.Range("J81").Value = 78000
$ws.Range("K81").Value = 176100
$ws.Range("L81").Value = 104500
$ws.Range("M81").Value = 208900
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("D89").Value = 361500
$ws.Range("E89").Value = 264400
$ws.Range("F89").Value = 210300
$ws.Range("G89").Value = 245100
$ws.Range("H89").Value = 350700
$ws.Range("I89").Value = 315100
$ws.Range("J89").Value = 49300
$ws.Range("K89").Value = 225200
$ws.Range("L89").Value = 450000
$ws.Range("M89").Value = 119300
$ws.Range("D91").Value = -60500
$ws.Range("E91").Value = -47000
$ws.Range("F91").Value = -223800
$ws.Range("G91").Value = -110600
$ws.Range("H91").Value = -104800
$ws.Range("I91").Value = -84800
$ws.Range("J91").Value = -74800
$ws.Range("K91").Value = -130200
$ws.Range("L91").Value = -227400
$ws.Range("M91").Value = -1400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("D94").Value = -75100
$ws.Range("E94").Value = 44500
$ws.Range("F94").Value = -1209100
$ws.Range("G94").Value = -1526200
$ws.Range("H94").Value = -103500
$ws.Range("I94").Value = -84500
$ws.Range("J94").Value = -67500
$ws.Range("K94").Value = 41200
$ws.Range("L94").Value = -227900
$ws.Range("M94").Value = 204500
$ws.Range("D96").Value = -3400
$ws.Range("E96").Value = -100
$ws.Range("F96").Value = -255800
$ws.Range("G96").Value = -80700
$ws.Range("H96").Value = -2000
$ws.Range("I96").Value = -6100
$ws.Range("J96").Value = -284100
$ws.Range("K96").Value = -93700
$ws.Range("L96").Value = -2700
$ws.Range("M96").Value = -1600
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = -233600
$ws.Range("E100").Value = -157400
$ws.Range("F100").Value = 540000
$ws.Range("G100").Value = 1271900
$ws.Range("H100").Value = -30600
$ws.Range("I100").Value = -19000
$ws.Range("J100").Value = -313400
$ws.Range("K100").Value = -106700
$ws.Range("L100").Value = -182500
$ws.Range("M100").Value = -269300
$ws.Range("D101").Value = 14500
$ws.Range("E101").Value = -8800
$ws.Range("F101").Value = 3400
$ws.Range("G101").Value = -1500
$ws.Range("H101").Value = 100
$ws.Range("I101").Value = -1500
$ws.Range("J101").Value = -1800
$ws.Range("K101").Value = 6200
$ws.Range("L101").Value = -700
$ws.Range("M101").Value = -2800
$ws.Range("D102").Value = 67300
$ws.Range("E102").Value = 142700
$ws.Range("F102").Value = -455400
$ws.Range("G102").Value = -10700
$ws.Range("H102").Value = 216600
$ws.Range("I102").Value = 210100
$ws.Range("J102").Value = -333300
$ws.Range("K102").Value = 165900
$ws.Range("L102").Value = 38800
$ws.Range("M102").Value = 51800
